# feat: add 2022-Q4 data
#
# 1. Insert a new "2022-Q4" worksheet (fund holdings detail) right after
#    the "总计" (summary) sheet and before the existing "2022-Q1" sheet.
# 2. Populate it with the 2022-Q4 fund holdings table.
# 3. Insert a new row into the "总计" sheet for the 2022-Q4 quarter,
#    pushing the existing 2022-Q1 / 2021-Q1 summary rows down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q4" sheet positioned before "2022-Q1"
# ---------------------------------------------------------------------
$q1Sheet = $wb.Worksheets.Item("2022-Q1")
$newSheet = $wb.Worksheets.Add($q1Sheet)
$newSheet.Name = "2022-Q4"

# Pick up the look (bold/border/centered) of the existing fund sheets'
# header row and index column so the new sheet matches its siblings.
$q1Sheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$q1Sheet.Range("A2:A3").Copy()
$newSheet.Range("A2:A3").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# ---------------------------------------------------------------------
# 2) Fill in the 2022-Q4 fund holdings data
# ---------------------------------------------------------------------
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Columns B:G hold plain text in this workbook (fund codes / formatted
# numbers kept verbatim, e.g. "005189", "1.00"), so force text formatting
# before writing them, otherwise Excel would coerce these into numbers.
$newSheet.Range("B2:G3").NumberFormat = "@"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "005189"
$newSheet.Range("C2").Value = "海富通量化前锋股票A"
$newSheet.Range("D2").Value = "0.54"
$newSheet.Range("E2").Value = "88.13"
$newSheet.Range("F2").Value = "1.00"
$newSheet.Range("G2").Value = "0.0054"
$newSheet.Range("H2").Value = 6

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "005188"
$newSheet.Range("C3").Value = "海富通量化前锋股票C"
$newSheet.Range("D3").Value = "0.03"
$newSheet.Range("E3").Value = "88.13"
$newSheet.Range("F3").Value = "1.00"
$newSheet.Range("G3").Value = "0.0003"
$newSheet.Range("H3").Value = 6

# ---------------------------------------------------------------------
# 3) Update the "总计" (summary) sheet: shift 2022-Q1 / 2021-Q1 rows down
#    one row and insert the new 2022-Q4 summary row at row 2.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

# Give the new row 4 (A4) the same index-column look as rows 2/3 (style
# "s=2" -> bold, bordered, centered), since it did not exist before.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# 2021-Q1 (was row 3) -> row 4
$totalSheet.Range("B4").Value = $totalSheet.Range("B3").Value()
$totalSheet.Range("C4").Value = $totalSheet.Range("C3").Value()
$totalSheet.Range("D4").Value = $totalSheet.Range("D3").Value()
$totalSheet.Range("A4").Value = 2

# 2022-Q1 (was row 2) -> row 3
$totalSheet.Range("B3").Value = $totalSheet.Range("B2").Value()
$totalSheet.Range("C3").Value = $totalSheet.Range("C2").Value()
$totalSheet.Range("D3").Value = $totalSheet.Range("D2").Value()
$totalSheet.Range("A3").Value = 1

# 2022-Q4 (new) -> row 2
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.01
